# Update the ticker lists on Sheet1 per the target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B values (row 2 .. row 23)
$colB = @(
    "NSE:63MOONS",
    "NSE:BEML",
    "NSE:CAMLINFINE",
    "NSE:CHOLAHLDNG",
    "NSE:CSLFINANCE",
    "NSE:DEEPINDS",
    "NSE:GNA",
    "NSE:HARDWYN",
    "NSE:IIFLSEC",
    "NSE:JAYSREETEA",
    "NSE:JKCEMENT",
    "NSE:JUBLINGREA",
    "NSE:KAJARIACER",
    "NSE:KIMS",
    "NSE:KRISHANA",
    "NSE:MARKSANS",
    "NSE:POKARNA",
    "NSE:PRSMJOHNSN",
    "NSE:RAMASTEEL",
    "NSE:RBL",
    "NSE:REFEX",
    "NSE:RGL"
)

# Column C values (row 2 .. row 23); empty string means no text (blank inlineStr)
$colC = @(
    "NSE:ARTEMISMED",
    "NSE:BRNL",
    "NSE:CANTABIL",
    "NSE:CHEVIOT",
    "NSE:CYIENTDLM",
    "NSE:GABRIEL",
    "NSE:GOACARBON",
    "NSE:GREAVESCOT",
    "NSE:HITECH",
    "NSE:JTLIND",
    "NSE:MAHLOG",
    "NSE:MODIRUBBER",
    "NSE:ONWARDTEC",
    "NSE:OPTIEMUS",
    "NSE:OSWALSEEDS",
    "NSE:PRICOLLTD",
    "NSE:PROZONER",
    "NSE:RITCO",
    "NSE:RSYSTEMS",
    "NSE:RVNL",
    "",
    ""
)

# Column E values (row 2 .. row 23); only rows 2 and 3 get new values
$colE = @(
    "NSE:ONGC",
    "NSE:PNB",
    "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", ""
)

# Column F values (row 2 .. row 23); only row 2 changes
$colF = @(
    "NSE:PIDILITIND",
    "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", "", ""
)

for ($i = 0; $i -lt $colB.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $colB[$i]
    if ($colC[$i] -ne "") {
        $ws.Cells.Item($row, 3).Value = $colC[$i]
    }
    if ($colE[$i] -ne "") {
        $ws.Cells.Item($row, 5).Value = $colE[$i]
    }
    if ($colF[$i] -ne "") {
        $ws.Cells.Item($row, 6).Value = $colF[$i]
    }
}

# New rows 22 and 23 need column A populated too (0-based sequence continues: 20, 21).
# Copy formatting from the last existing numbered row (A21) so the new cells
# reuse the same cell style (bold, centered, bordered) instead of minting a
# brand new style entry.
$ws.Range("A21").Copy($ws.Range("A22"))
$ws.Range("A21").Copy($ws.Range("A23"))
$ws.Cells.Item(22, 1).Value = 20
$ws.Cells.Item(23, 1).Value = 21
